# This script applies the latest cryptocurrency price/volume scrape to
# the "cryptos" worksheet (GitHub Actions scheduled update).
#
# Every touched cell is plain text (coin names, URLs, price strings such
# as "65.975.56", and padded percent strings such as "  +7.12%  "). Some
# of the new price strings parse as plain numbers (e.g. "1.00", "2.40"),
# so before writing we force the cell format to Text ("@") to stop Excel
# from silently reinterpreting/rounding them, then restore the default
# "Normal" style afterwards so no stray formatting is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $rng = $ws.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "65.975.56"
Set-TextValue "E2" "  +7.12%  "
# Row 3
Set-TextValue "D3" "3.010.53"
Set-TextValue "E3" "  +3.97%  "
# Row 4
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.16%  "
# Row 5
Set-TextValue "D5" "586.29"
Set-TextValue "E5" "  +3.09%  "
# Row 6
Set-TextValue "D6" "154.31"
Set-TextValue "E6" "  +6.78%  "
# Row 7
Set-TextValue "E7" "  -0.13%  "
# Row 8
Set-TextValue "D8" "3.003.55"
Set-TextValue "E8" "  +3.81%  "
# Row 9
Set-TextValue "E9" "  +2.20%  "
# Row 10
Set-TextValue "D10" "6.98"
Set-TextValue "E10" "  +0.27%  "
# Row 11
Set-TextValue "E11" "  +4.23%  "
# Row 12
Set-TextValue "D12" "0.449"
Set-TextValue "E12" "  +3.76%  "
# Row 13
Set-TextValue "E13" "  +2.98%  "
# Row 14
Set-TextValue "E14" "  +6.35%  "
# Row 15
Set-TextValue "E15" "  +0.74%  "
# Row 16
Set-TextValue "D16" "65.921.35"
Set-TextValue "E16" "  +7.07%  "
# Row 17
Set-TextValue "D17" "3.507.13"
Set-TextValue "E17" "  +3.88%  "
# Row 18
Set-TextValue "E18" "  +5.74%  "
# Row 19
Set-TextValue "D19" "3.009.08"
Set-TextValue "E19" "  +3.86%  "
# Row 20
Set-TextValue "D20" "457.96"
Set-TextValue "E20" "  +5.73%  "
# Row 21
Set-TextValue "E21" "  +4.75%  "
# Row 22
Set-TextValue "E22" "  +3.94%  "
# Row 23
Set-TextValue "E23" "  +7.19%  "
# Row 24
Set-TextValue "D24" "81.62"
Set-TextValue "E24" "  +2.83%  "
# Row 25
Set-TextValue "D25" "12.54"
Set-TextValue "E25" "  +3.80%  "
# Row 26
Set-TextValue "D26" "2.25"
Set-TextValue "E26" "  +10.76%  "
# Row 27
Set-TextValue "D27" "10.69"
Set-TextValue "E27" "  +6.45%  "
# Row 29
Set-TextValue "D29" "2.40"
Set-TextValue "E29" "  +16.52%  "
# Row 30
Set-TextValue "E30" "  +10.98%  "
# Row 31
Set-TextValue "D31" "2.61"
Set-TextValue "E31" "  +3.96%  "
# Row 32
Set-TextValue "E32" "  -3.01%  "
# Row 33
Set-TextValue "E33" "  +5.88%  "
# Row 34
Set-TextValue "E34" "  +3.87%  "
# Row 35
Set-TextValue "D35" "0.999"
Set-TextValue "E35" "  -0.03%  "
# Row 36
Set-TextValue "D36" "0.993"
Set-TextValue "E36" "  +3.32%  "
# Row 37
Set-TextValue "E37" "  +6.99%  "
# Row 38
Set-TextValue "D38" "2.14"
Set-TextValue "E38" "  +10.23%  "
# Row 39
Set-TextValue "D39" "45.55"
Set-TextValue "E39" "  +14.43%  "
# Row 40
Set-TextValue "E40" "  +0.70%  "
# Row 41
Set-TextValue "E41" "  +3.04%  "
# Row 42
Set-TextValue "D42" "0.122"
Set-TextValue "E42" "  +5.91%  "
# Row 43
Set-TextValue "E43" "  +12.11%  "
# Row 44
Set-TextValue "D44" "8.45"
Set-TextValue "E44" "  +2.63%  "
# Row 45
Set-TextValue "D45" "387.58"
Set-TextValue "E45" "  +12.14%  "
# Row 46
Set-TextValue "B46" "Maker"
Set-TextValue "C46" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D46" "2.766.75"
Set-TextValue "E46" "  +2.05%  "
# Row 47
Set-TextValue "B47" "VeChain"
Set-TextValue "C47" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D47" "0.0352"
Set-TextValue "E47" "  +4.55%  "
# Row 48
Set-TextValue "D48" "135.43"
Set-TextValue "E48" "  +1.90%  "
# Row 50
Set-TextValue "E50" "  +8.16%  "
# Row 51
Set-TextValue "E51" "  +2.87%  "
